$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: "Average of SW(S*)/SW(OPT)" | AVERAGE(N2:N11), bold sz12, vertical-center ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$f = $ws.Range("B14").Font
$f.Bold = $true
$f.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# --- Row 12: J12 = AVERAGE(J2:J11), bold ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Reuse the exact same cell format as B14 for the rows below (avoids creating redundant style entries)
$ws.Range("B14").Copy()

# --- Row 15: "Average of SC(S*)/SC(OPT)" | AVERAGE(Z2:Z11) ---
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# --- Row 16: "Worst of SW(S*)/SW(OPT)" | MIN(N2:N11) ---
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# --- Row 17: "Worst of SC(S*)/SC(OPT)" | MAX(Z2:Z11) ---
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the active selection on J12, matching the saved view state ---
$excel.CutCopyMode = $false
[void]$ws.Range("J12").Select()
